$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PaymentPage")

# Change A3 and A4 from numeric values to text values "6" and "9"
$ws.Range("A3").Value = "6"
$ws.Range("A4").Value = "9"

# Update selection / active cell on the sheet
$ws.Range("A1:C4").Select()
$ws.Range("A4").Activate()

# Adjust the window view size/position
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Left = 2340
$wb.Windows.Item(1).Top = 2820
$wb.Windows.Item(1).Width = 15375
$wb.Windows.Item(1).Height = 7965
